# Update "想去人数" (interest count, column F) values on the 展览 and
# 全部类型 worksheets to reflect the newly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 3014
$ws1.Range("F5").Value  = 460
$ws1.Range("F7").Value  = 26
$ws1.Range("F8").Value  = 259
$ws1.Range("F10").Value = 14376
$ws1.Range("F13").Value = 5749
$ws1.Range("F17").Value = 59
$ws1.Range("F25").Value = 10539

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 3014
$ws4.Range("F6").Value  = 460
$ws4.Range("F8").Value  = 26
$ws4.Range("F9").Value  = 259
$ws4.Range("F11").Value = 14376
$ws4.Range("F14").Value = 5749
$ws4.Range("F18").Value = 59
$ws4.Range("F27").Value = 10539
